# EJERCICIO 12_UNIDAD 1.xlsx
# "problemas resueltos hasta el 14"
#
# The workbook previously pulled its summary figures (J2:J5, M2:M7, I7, I8)
# from an external workbook via SUMIF/AVERAGE/SUM formulas that referenced
# '[1].xlsx]EJERCICIO-12'. That external link is now replaced with formulas
# that work off this sheet's own data (columns E/F for payment-method
# totals, B/F for product totals), and the now-unused external link is
# broken/removed. One more bug was introduced directly into F16 (typed over
# the shared formula with a typo), and the sheet view was re-zoomed/re-
# selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recaudación por medio de pago (H1:J8) -------------------------------
$ws.Range("J2").Formula = '=SUMIF(E10:E34,"MERCADO PAGO",F10:F34)'
$ws.Range("J3").Formula = '=SUMIF(E10:E34,"TARJETA DE CRÉDITO",F10:F34)'
$ws.Range("J4").Formula = '=SUMIF(E10:E34,"TARJETA DE DEBITO",F10:F34)'
$ws.Range("J5").Formula = '=SUMIF(E10:E34,"EFECTIVO",F10:F34)'

# Precio promedio (now resolves instead of #DIV/0!, since M2:M7 are local)
$ws.Range("I7").Formula = '=AVERAGE(M2:M7)'

# Total recaudado — now errors (#VALUE!) because F29 holds a #VALUE! error
$ws.Range("I8").Formula = '=SUM(F10:F34)'

# --- Recaudación por producto (L1:M7) ------------------------------------
# Style 12 ("$" #,##0.00) -> style 8 ([$$-2C0A] #,##0.00, same as column B)
$moneyFmt = '[$$-2C0A]\ #,##0.00'

$ws.Range("M2").Formula = '=SUMIF($B$2:$B$51,"Pan lactal",$F$10:$F$59)'
$ws.Range("M2").NumberFormat = $moneyFmt

$ws.Range("M3").Formula = '=SUMIF($B$2:$B$51,"Galletas",$F$10:$F$59)'
$ws.Range("M3").NumberFormat = $moneyFmt

$ws.Range("M4").Formula = '=SUMIF($B$2:$B$51,"Pan ralladoo",$F$10:$F$59)'
$ws.Range("M4").NumberFormat = $moneyFmt

$ws.Range("M5").Formula = '=SUMIF($B$2:$B$51,"Budín marmolado",$F$10:$F$59)'
$ws.Range("M5").NumberFormat = $moneyFmt

$ws.Range("M6").Formula = '=SUMIF($B$2:$B$51,"Pan de pebete",$F$10:$F$59)'
$ws.Range("M6").NumberFormat = $moneyFmt

$ws.Range("M7").Formula = '=SUMIF($B$2:$B$51,"Prepizza",$F$10:$F$59)'
$ws.Range("M7").NumberFormat = $moneyFmt

# --- F16: typed directly over the shared formula, with a typo (D16*D16
# instead of D16*B16), detaching it from the F11:F34 shared formula group.
$ws.Range("F16").Formula = '=IF(OR(E16="MERCADO PAGO",E16="TARJETA DE CRÉDITO"),(D16*C16)*$F$2,IF(E16="EFECTIVO",(D16*C16)*$E$3,D16*D16))'

# --- Drop the now-unused external workbook link --------------------------
# (removes xl/externalLinks/externalLink1.xml, its rels, the
# <externalReferences> block in xl/workbook.xml, and the Content_Types
# override) now that no formula references '[1].xlsx]...' any more.
$wb.BreakLink("file:///C:\Users\Alan\AppData\Local\Microsoft\Windows\INetCache\IE\PYYAOL85\Ejercicios-nivelación%5b1%5d.xlsx")

# --- View state: selection moved to J4, zoomed in to 160% ----------------
$ws.Range("J4").Select()
$excel.ActiveWindow.Zoom = 160
